$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("mesa" 1): drop the stray empty Usuario entry, table is free again
$ws.Range("B5").ClearContents()

# New row 6: mesa 6, Libre, capacity 12
$ws.Range("A6").Value = 6
$ws.Range("C6").Value = "Libre"
$ws.Range("D6").Value = "'12"
$ws.Range("D6").ClearFormats()

# New row 7: mesa 7, Libre, capacity 3
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "'"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = "Libre"
$ws.Range("D7").Value = "'3"
$ws.Range("D7").ClearFormats()
